$wb = $excel.ActiveWorkbook

# Values for F2:F9 ("想去人数" - number of people wanting to go)
$values = @(348, 90, 1530, 21, 47, 129, 53, 347)

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 6).Value = $values[$i]
    }
}
